$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new cell N2, formatted like D2/M2 (border-only style), no value ---
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

# --- Row 3: new cell N3 = 2021, formatted like M3 (year header style) ---
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2021

# --- Row 4: bold the whole data row (D4:L4 pick up the bold variant of their
#     existing style, matching M4 which was already bold), then add N4 ---
$ws.Range("D4:M4").Font.Bold = $true
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 95.134712433469176

# --- Rows 5-14: new N cell formatted like the D column of that row ---
$ws.Range("D5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 99.705541665880986

$ws.Range("D6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 92.425193326577897

$ws.Range("D7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 88.209991167538519

$ws.Range("D8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 92.225038985690773

$ws.Range("D9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 96.801032063987265

$ws.Range("D10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 97.660491031729507

$ws.Range("D11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 90.23262877800066

$ws.Range("D12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 99.653994395099105

$ws.Range("D13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 100

$ws.Range("D14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 100

# --- Row 15: new N15 formatted like M15/L15 (bottom border + right align) ---
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = 100

# --- Selection moves to the newly added N2 cell ---
$ws.Range("N2").Select()

# --- Print vertical resolution bumped from 0 (unset/96dpi) to 300 ---
$ws.PageSetup.PrintQuality = 300
